# DESIGN/rules/sp1/Main.xlsx - project sp1 save.
# Business edit: "Integer min" rule value (C8) updated from 0 to 1110.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1110

# Leave the selection where the author ended up before saving.
$ws.Range("D14").Select()
